# Auto-generated edit script: appends 6 new daily rows (2026-01-06 .. 2026-01-08)
# for the two charging stations, mirroring the existing row layout/styles.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 256
$ws.Range("A256").Value = 46028
$ws.Range("B256").Value = "四方坪站充电量(kw)"
$ws.Range("C256").Value = 905.92
$ws.Range("D256").Value = 1399.95
$ws.Range("E256").Value = 596.28
$ws.Range("F256").Value = 641.79999999999995
$ws.Range("G256").Value = 428.53
$ws.Range("H256").Value = 751.2
$ws.Range("I256").Value = 542.78
$ws.Range("J256").Value = 183.87
$ws.Range("K256").Value = 261.31
$ws.Range("L256").Value = 148.96
$ws.Range("M256").Value = 176.34
$ws.Range("N256").Value = 432.99
$ws.Range("O256").Value = 1121.33
$ws.Range("P256").Value = 1678.17
$ws.Range("Q256").Value = 686.2
$ws.Range("R256").Value = 663.69
$ws.Range("S256").Value = 387.36
$ws.Range("T256").Value = 143.19999999999999
$ws.Range("U256").Value = 172.85
$ws.Range("V256").Value = 115.21
$ws.Range("W256").Value = 237.71
$ws.Range("X256").Value = 93.87
$ws.Range("Y256").Value = 79.38
$ws.Range("Z256").Value = 19.86

# Row 257
$ws.Range("A257").Value = 46028
$ws.Range("B257").Value = "高岭站充电量(kw)"
$ws.Range("C257").Value = 690.1099999999999
$ws.Range("D257").Value = 166.51
$ws.Range("E257").Value = 115.38
$ws.Range("F257").Value = 4.96
$ws.Range("G257").Value = 84.449999999999989
$ws.Range("H257").Value = 118.50999999999999
$ws.Range("I257").Value = 80.819999999999993
$ws.Range("J257").Value = 19.37
$ws.Range("K257").Value = 143.24
$ws.Range("L257").Value = 175.35
$ws.Range("M257").Value = 299.40999999999997
$ws.Range("N257").Value = 398.64000000000004
$ws.Range("O257").Value = 534.76
$ws.Range("P257").Value = 426.53999999999996
$ws.Range("Q257").Value = 396.89
$ws.Range("R257").Value = 152.45000000000002
$ws.Range("S257").Value = 247.23
$ws.Range("T257").Value = 74.27
$ws.Range("U257").Value = 89.48
$ws.Range("V257").Value = 237.32999999999998
$ws.Range("W257").Value = 19.809999999999999
$ws.Range("X257").Value = 0
$ws.Range("Y257").Value = 63.54
$ws.Range("Z257").Value = 74.13

# Row 258
$ws.Range("A258").Value = 46029
$ws.Range("B258").Value = "四方坪站充电量(kw)"
$ws.Range("C258").Value = 839.28
$ws.Range("D258").Value = 1306.94
$ws.Range("E258").Value = 641.28
$ws.Range("F258").Value = 883.64
$ws.Range("G258").Value = 517.36
$ws.Range("H258").Value = 731.45
$ws.Range("I258").Value = 739.43
$ws.Range("J258").Value = 172.45
$ws.Range("K258").Value = 249
$ws.Range("L258").Value = 107.66
$ws.Range("M258").Value = 118.81
$ws.Range("N258").Value = 348.56
$ws.Range("O258").Value = 871.67
$ws.Range("P258").Value = 1576.45
$ws.Range("Q258").Value = 677.31
$ws.Range("R258").Value = 580.20000000000005
$ws.Range("S258").Value = 142.65
$ws.Range("T258").Value = 243.89
$ws.Range("U258").Value = 208.96
$ws.Range("V258").Value = 94.17
$ws.Range("W258").Value = 193.16
$ws.Range("X258").Value = 107.43
$ws.Range("Y258").Value = 64.81
$ws.Range("Z258").Value = 24.85

# Row 259
$ws.Range("A259").Value = 46029
$ws.Range("B259").Value = "高岭站充电量(kw)"
$ws.Range("C259").Value = 617.50000000000011
$ws.Range("D259").Value = 92.609999999999985
$ws.Range("E259").Value = 16.59
$ws.Range("F259").Value = 0
$ws.Range("G259").Value = 52.510000000000005
$ws.Range("H259").Value = 157.5
$ws.Range("I259").Value = 60.33
$ws.Range("J259").Value = 30.61
$ws.Range("K259").Value = 221.05000000000004
$ws.Range("L259").Value = 486.21
$ws.Range("M259").Value = 102.57
$ws.Range("N259").Value = 198.82
$ws.Range("O259").Value = 658.37999999999977
$ws.Range("P259").Value = 426.19
$ws.Range("Q259").Value = 392.35
$ws.Range("R259").Value = 138.57
$ws.Range("S259").Value = 178.69
$ws.Range("T259").Value = 131.47000000000003
$ws.Range("U259").Value = 46.1
$ws.Range("V259").Value = 132.35
$ws.Range("W259").Value = 99.889999999999986
$ws.Range("X259").Value = 104.58
$ws.Range("Y259").Value = 40.510000000000005
$ws.Range("Z259").Value = 47.41

# Row 260
$ws.Range("A260").Value = 46030
$ws.Range("B260").Value = "四方坪站充电量(kw)"
$ws.Range("C260").Value = 849.21
$ws.Range("D260").Value = 1979.69
$ws.Range("E260").Value = 722.08
$ws.Range("F260").Value = 266.36
$ws.Range("G260").Value = 426.46
$ws.Range("H260").Value = 790.76
$ws.Range("I260").Value = 419.13
$ws.Range("J260").Value = 58.8
$ws.Range("K260").Value = 212.71
$ws.Range("L260").Value = 86.58
$ws.Range("M260").Value = 185.32
$ws.Range("N260").Value = 332.63
$ws.Range("O260").Value = 681.36
$ws.Range("P260").Value = 2203.21
$ws.Range("Q260").Value = 688.89
$ws.Range("R260").Value = 307.3
$ws.Range("S260").Value = 156.75
$ws.Range("T260").Value = 191.22
$ws.Range("U260").Value = 220.49
$ws.Range("V260").Value = 100.48
$ws.Range("W260").Value = 91.96
$ws.Range("X260").Value = 24.72
$ws.Range("Y260").Value = 30.56
$ws.Range("Z260").Value = 51.66

# Row 261
$ws.Range("A261").Value = 46030
$ws.Range("B261").Value = "高岭站充电量(kw)"
$ws.Range("C261").Value = 558.13
$ws.Range("D261").Value = 183.7
$ws.Range("E261").Value = 98.080000000000013
$ws.Range("F261").Value = 128.32999999999998
$ws.Range("G261").Value = 77.17
$ws.Range("H261").Value = 151.72999999999999
$ws.Range("I261").Value = 71.37
$ws.Range("J261").Value = 83.54
$ws.Range("K261").Value = 162.66
$ws.Range("L261").Value = 107.35
$ws.Range("M261").Value = 186.85
$ws.Range("N261").Value = 327.39
$ws.Range("O261").Value = 454.59
$ws.Range("P261").Value = 473.75000000000011
$ws.Range("Q261").Value = 167.35000000000002
$ws.Range("R261").Value = 207.33
$ws.Range("S261").Value = 101.39
$ws.Range("T261").Value = 227.5
$ws.Range("U261").Value = 0
$ws.Range("V261").Value = 31.64
$ws.Range("W261").Value = 129.97
$ws.Range("X261").Value = 62.97
$ws.Range("Y261").Value = 14.28
$ws.Range("Z261").Value = 5.51

# Update the view state to match the scrolled/selected position after the append
$ws.Range("F264").Select() | Out-Null
